$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C4"  = -12.481
    "C7"  = -12.793
    "D7"  = -7.7
    "D15" = -8.24
    "C16" = -13.383
    "D21" = -8.1
    "D22" = -7.900000000000001
    "D23" = -7.856
    "C28" = -13.219
    "C29" = -11.997
    "C32" = -13.434
    "D34" = -7.782999999999999
    "C40" = -12.232
    "D43" = -7.707000000000001
    "D45" = -7.498
    "D50" = -8.178000000000001
    "D51" = -8.370000000000001
    "C52" = -11.614
    "C57" = -13.622
    "C66" = -11.127
    "D66" = -7.447
    "D67" = -7.141999999999999
    "D79" = -7.532000000000001
    "D84" = -8.141
    "D92" = -6.592999999999999
    "D97" = -8.313000000000001
    "C100" = -12.868
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$wb.Save()
